$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Turn the per-row L/M formulas into shared formulas (same formula text, now
#     written once across the whole range so Excel stores them as t="shared") ---
$ws.Range("L2:L14").Formula = "=110-25*((E2/C2)/(F2/D2))"
$ws.Range("M2:M14").Formula = "=ABS(L2-G2)"

# --- Drop the old scratch column F (rows 30:42) that just duplicated column G ---
$ws.Range("F30:F42").ClearContents()

# --- New "Model prediction" regression summary block (rows 24-26) ---
$ws.Range("C24").Value = "Coefficient"
$ws.Range("C25").Value = "RDC"
$ws.Range("D25").Value = "IRDC"
$ws.Range("E25").Value = "RACrms"
$ws.Range("F25").Value = "IRACrms"
$ws.Range("G25").Value = "Interception"
$ws.Range("H25").Value = "Corrélation"

$ws.Range("C26").Value = 44.101038809999999
$ws.Range("D26").Value = -24.604666720000001
$ws.Range("E26").Value = -4575.9874201800003
$ws.Range("F26").Value = 2380.0554081800001
$ws.Range("G26").Value = 80.280442275973797
$ws.Range("H26").Value = 0.99422622728771903

# Header style: centered, merged across C24:H24
$ws.Range("C24:H24").HorizontalAlignment = -4108
$ws.Range("C24:H24").Merge()

# --- Column widths tweaked by the author ---
$ws.Columns("B").ColumnWidth = 11.5
$ws.Columns("G").ColumnWidth = 10.33

# --- Selection left on the new summary table ---
$null = $ws.Range("C24:H26").Select()
